# Regenerate save_data: column G ("K") used to hold a called-strikes
# count ("Strike#"); it is recomputed here as the strikeout ("K") count
# per outing and written back row by row (rows 2-76, sheet "Sheet1").
#
# s_vals: row number (1-based) -> recomputed K value
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sVals = [ordered]@{
    2  = 1
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    8  = 1
    9  = 2
    10 = 0
    11 = 3
    12 = 2
    13 = 1
    14 = 1
    15 = 1
    16 = 1
    17 = 1
    18 = 2
    19 = 2
    20 = 1
    21 = 0
    22 = 0
    23 = 1
    24 = 0
    25 = 0
    26 = 1
    27 = 1
    28 = 1
    29 = 0
    30 = 0
    31 = 1
    32 = 0
    33 = 0
    34 = 2
    35 = 0
    36 = 0
    37 = 2
    38 = 1
    39 = 0
    40 = 0
    41 = 1
    42 = 1
    43 = 0
    44 = 1
    45 = 0
    46 = 1
    47 = 0
    48 = 1
    49 = 2
    50 = 3
    51 = 1
    52 = 0
    53 = 1
    54 = 0
    55 = 1
    56 = 1
    57 = 1
    58 = 2
    59 = 1
    60 = 0
    61 = 1
    62 = 1
    63 = 0
    64 = 1
    65 = 1
    66 = 0
    67 = 0
    68 = 1
    69 = 1
    70 = 0
    71 = 3
    72 = 1
    73 = 0
    74 = 1
    75 = 1
    76 = 1
}

foreach ($row in $sVals.Keys) {
    $ws.Cells.Item($row, 7).Value = $sVals[$row]
}
